$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.840815544128418
$ws.Range("B1").Value = 1.498015880584717
$ws.Range("C1").Value = 6.170209884643555
$ws.Range("D1").Value = 2.968657970428467
$ws.Range("E1").Value = 1.612660765647888
